$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-12-11 01:57:14"

for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}
